$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11  = -21.74759999999999
    12  = -21.402
    15  = -21.65270000000001
    27  = -21.87449999999999
    28  = -21.9481
    31  = -21.73320000000001
    32  = -21.29839999999998
    36  = -20.32249999999999
    38  = -19.98009999999999
    46  = -21.8788
    54  = -21.83759999999999
    55  = -22.02780000000001
    56  = -22.0272
    67  = -21.49089999999997
    69  = -21.66329999999998
    72  = -21.6962
    73  = -20.11509999999999
    83  = -21.58209999999998
    86  = -21.8572
    91  = -20.65999999999998
    93  = -21.4556
    99  = -21.77220000000001
    104 = -21.26339999999999
    105 = -19.87809999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
